$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value2 = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "29.915.95"
Set-TextValue "E2" "  +0.40%  "
Set-TextValue "D3" "1.888.47"
Set-TextValue "E3" "  -0.01%  "
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "0.7746"
Set-TextValue "E5" "  +0.65%  "
Set-TextValue "D6" "242.91"
Set-TextValue "E6" "  -0.55%  "
Set-TextValue "E7" "  +0.06%  "
Set-TextValue "D8" "0.3104"
Set-TextValue "E8" "  -0.63%  "
Set-TextValue "D9" "25.64"
Set-TextValue "E9" "  +1.63%  "
Set-TextValue "D10" "0.07143"
Set-TextValue "E10" "  -0.82%  "
Set-TextValue "D11" "0.08546"
Set-TextValue "E11" "  +5.43%  "
Set-TextValue "D12" "0.7643"
Set-TextValue "E12" "  -0.14%  "
Set-TextValue "D13" "1.943.04"
Set-TextValue "E13" "  +1.13%  "
Set-TextValue "D14" "5.349"
Set-TextValue "E14" "  -2.92%  "
Set-TextValue "D15" "93.79"
Set-TextValue "E15" "  +1.68%  "
Set-TextValue "E16" "  +0.26%  "
Set-TextValue "D17" "29.988.57"
Set-TextValue "E17" "  +0.61%  "
Set-TextValue "E18" "  -1.07%  "
Set-TextValue "D19" "244.16"
Set-TextValue "E19" "  +0.54%  "
Set-TextValue "B20" "ShibaInu"
Set-TextValue "C20" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D20" "0.000007807"
Set-TextValue "E20" "  +0.53%  "
Set-TextValue "B21" "WrappedliquidstakedEther2.0"
Set-TextValue "C21" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D21" "2.232.36"
Set-TextValue "E21" "  +4.21%  "
Set-TextValue "D22" "0.9981"
Set-TextValue "E22" "  -0.23%  "
Set-TextValue "D23" "7.926"
Set-TextValue "E23" "  -3.09%  "
Set-TextValue "D24" "1.001"
Set-TextValue "E24" "  +0.01%  "
Set-TextValue "D25" "0.1637"
Set-TextValue "E25" "  +4.85%  "
Set-TextValue "D26" "9.359"
Set-TextValue "E26" "  -0.48%  "
Set-TextValue "D27" "162.39"
Set-TextValue "E27" "  +0.12%  "
Set-TextValue "E28" "  +0.25%  "
Set-TextValue "D29" "2.036"
Set-TextValue "E29" "  -0.16%  "
Set-TextValue "E30" "  -1.20%  "
Set-TextValue "D31" "1.536"
Set-TextValue "E31" "  -0.86%  "
Set-TextValue "D32" "4.505"
Set-TextValue "E32" "  +1.22%  "
Set-TextValue "E33" "  +0.67%  "
Set-TextValue "D34" "0.05438"
Set-TextValue "E34" "  -1.78%  "
Set-TextValue "E35" "  -1.22%  "
Set-TextValue "D36" "0.7459"
Set-TextValue "E36" "  -0.21%  "
Set-TextValue "D37" "1.003"
Set-TextValue "E37" "  +0.33%  "
Set-TextValue "D38" "2.694"
Set-TextValue "E38" "  +2.40%  "
Set-TextValue "D39" "0.01963"
Set-TextValue "E39" "  +2.28%  "
Set-TextValue "D40" "2.784"
Set-TextValue "E40" "  +0.16%  "
Set-TextValue "D41" "0.4464"
Set-TextValue "E41" "  +1.16%  "
Set-TextValue "D42" "1.108.36"
Set-TextValue "E42" "  -4.26%  "
Set-TextValue "B43" "Aave"
Set-TextValue "C43" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D43" "73.18"
Set-TextValue "E43" "  -0.44%  "
Set-TextValue "B44" "FraxShare"
Set-TextValue "C44" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D44" "6.083"
Set-TextValue "E44" "  +3.03%  "
Set-TextValue "D45" "0.8494"
Set-TextValue "E45" "  +0.01%  "
Set-TextValue "E46" "  +0.09%  "
Set-TextValue "D47" "103.65"
Set-TextValue "E47" "  +0.83%  "
Set-TextValue "B48" "RenderToken"
Set-TextValue "C48" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D48" "1.872"
Set-TextValue "E48" "  -0.50%  "
Set-TextValue "B49" "RocketPoolETH"
Set-TextValue "C49" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D49" "2.139.72"
Set-TextValue "E49" "  +4.95%  "
Set-TextValue "D50" "7.606"
Set-TextValue "E50" "  +2.24%  "
Set-TextValue "E51" "  -0.92%  "
